$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 118.0346986666667
$ws.Range("H2").Value2 = 354.104096
$ws.Range("I2").Value2 = 0.2666057129183408
$ws.Range("J2").Value2 = 0.2666057129183408
$ws.Range("M2").Value2 = 28.22141
$ws.Range("N2").Value2 = 84.66423
$ws.Range("O2").Value2 = 0.007710741921554872
$ws.Range("P2").Value2 = 0.007710741921554872
$ws.Range("Q2").Value2 = 3331.105625298453
$ws.Range("R2").Value2 = 29979.95062768608
$ws.Range("S2").Value2 = 0.002055727847125474
$ws.Range("T2").Value2 = 0.002055727847125474
$ws.Range("G3").Value2 = 118.0346986666667
$ws.Range("H3").Value2 = 354.104096
$ws.Range("I3").Value2 = 0.2666057129183408
$ws.Range("J3").Value2 = 0.2666057129183408
$ws.Range("O3").Value2 = 0.001803104886918205
$ws.Range("P3").Value2 = 0.001803104886918206
$ws.Range("Q3").Value2 = 778.9565378950186
$ws.Range("R3").Value2 = 7010.608841055167
$ws.Range("S3").Value2 = 0.0004807180638433725
$ws.Range("T3").Value2 = 0.0004807180638433725
$ws.Range("G4").Value2 = 118.0346986666667
$ws.Range("H4").Value2 = 354.104096
$ws.Range("I4").Value2 = 0.2666057129183408
$ws.Range("J4").Value2 = 0.2666057129183408
$ws.Range("M4").Value2 = 1.757996666666666
$ws.Range("N4").Value2 = 5.27399
$ws.Range("O4").Value2 = 0.000480325348578274
$ws.Range("P4").Value2 = 0.0004803253485782741
$ws.Range("Q4").Value2 = 207.5046068070044
$ws.Range("R4").Value2 = 1867.54146126304
$ws.Range("S4").Value2 = 0.0001280574819904613
$ws.Range("T4").Value2 = 0.0001280574819904613
$ws.Range("G5").Value2 = 118.0346986666667
$ws.Range("H5").Value2 = 354.104096
$ws.Range("I5").Value2 = 0.2666057129183408
$ws.Range("J5").Value2 = 0.2666057129183408
$ws.Range("M5").Value2 = 3623.433471666667
$ws.Range("N5").Value2 = 10870.300415
$ws.Range("O5").Value2 = 0.9900058278429487
$ws.Range("P5").Value2 = 0.9900058278429487
$ws.Range("Q5").Value2 = 427690.8779668888
$ws.Range("R5").Value2 = 3849217.901701999
$ws.Range("S5").Value2 = 0.2639412095253815
$ws.Range("T5").Value2 = 0.2639412095253815
$ws.Range("I6").Value2 = 0.4881754016778185
$ws.Range("J6").Value2 = 0.4881754016778186
$ws.Range("M6").Value2 = 28.22141
$ws.Range("N6").Value2 = 84.66423
$ws.Range("O6").Value2 = 0.007710741921554872
$ws.Range("P6").Value2 = 0.007710741921554872
$ws.Range("Q6").Value2 = 6099.508554639991
$ws.Range("R6").Value2 = 54895.57699175991
$ws.Range("S6").Value2 = 0.003764194534789044
$ws.Range("T6").Value2 = 0.003764194534789044
$ws.Range("I7").Value2 = 0.4881754016778185
$ws.Range("J7").Value2 = 0.4881754016778186
$ws.Range("O7").Value2 = 0.001803104886918205
$ws.Range("P7").Value2 = 0.001803104886918206
$ws.Range("S7").Value2 = 0.0008802314524385325
$ws.Range("T7").Value2 = 0.0008802314524385327
$ws.Range("I8").Value2 = 0.4881754016778185
$ws.Range("J8").Value2 = 0.4881754016778186
$ws.Range("M8").Value2 = 1.757996666666666
$ws.Range("N8").Value2 = 5.27399
$ws.Range("O8").Value2 = 0.000480325348578274
$ws.Range("P8").Value2 = 0.0004803253485782741
$ws.Range("Q8").Value2 = 379.9567671268699
$ws.Range("R8").Value2 = 3419.610904141829
$ws.Range("S8").Value2 = 0.0002344830199782371
$ws.Range("T8").Value2 = 0.0002344830199782372
$ws.Range("I9").Value2 = 0.4881754016778185
$ws.Range("J9").Value2 = 0.4881754016778186
$ws.Range("M9").Value2 = 3623.433471666667
$ws.Range("N9").Value2 = 10870.300415
$ws.Range("O9").Value2 = 0.9900058278429487
$ws.Range("P9").Value2 = 0.9900058278429487
$ws.Range("Q9").Value2 = 783134.6292619579
$ws.Range("R9").Value2 = 7048211.663357621
$ws.Range("S9").Value2 = 0.4832964926706127
$ws.Range("T9").Value2 = 0.4832964926706128
$ws.Range("G10").Value2 = 45.876452
$ws.Range("H10").Value2 = 137.629356
$ws.Range("I10").Value2 = 0.1036214293744632
$ws.Range("J10").Value2 = 0.1036214293744632
$ws.Range("M10").Value2 = 28.22141
$ws.Range("N10").Value2 = 84.66423
$ws.Range("O10").Value2 = 0.007710741921554872
$ws.Range("P10").Value2 = 0.007710741921554872
$ws.Range("Q10").Value2 = 1294.69816123732
$ws.Range("R10").Value2 = 11652.28345113588
$ws.Range("S10").Value2 = 0.0007989980994491107
$ws.Range("T10").Value2 = 0.0007989980994491108
$ws.Range("G11").Value2 = 45.876452
$ws.Range("H11").Value2 = 137.629356
$ws.Range("I11").Value2 = 0.1036214293744632
$ws.Range("J11").Value2 = 0.1036214293744632
$ws.Range("O11").Value2 = 0.001803104886918205
$ws.Range("P11").Value2 = 0.001803104886918206
$ws.Range("Q11").Value2 = 302.756415058472
$ws.Range("R11").Value2 = 2724.807735526248
$ws.Range("S11").Value2 = 0.0001868403056945442
$ws.Range("T11").Value2 = 0.0001868403056945443
$ws.Range("G12").Value2 = 45.876452
$ws.Range("H12").Value2 = 137.629356
$ws.Range("I12").Value2 = 0.1036214293744632
$ws.Range("J12").Value2 = 0.1036214293744632
$ws.Range("M12").Value2 = 1.757996666666666
$ws.Range("N12").Value2 = 5.27399
$ws.Range("O12").Value2 = 0.000480325348578274
$ws.Range("P12").Value2 = 0.0004803253485782741
$ws.Range("Q12").Value2 = 80.65064969449332
$ws.Range("R12").Value2 = 725.85584725044
$ws.Range("S12").Value2 = 0.00004977199918446803
$ws.Range("T12").Value2 = 0.00004977199918446804
$ws.Range("G13").Value2 = 45.876452
$ws.Range("H13").Value2 = 137.629356
$ws.Range("I13").Value2 = 0.1036214293744632
$ws.Range("J13").Value2 = 0.1036214293744632
$ws.Range("M13").Value2 = 3623.433471666667
$ws.Range("N13").Value2 = 10870.300415
$ws.Range("O13").Value2 = 0.9900058278429487
$ws.Range("P13").Value2 = 0.9900058278429487
$ws.Range("Q13").Value2 = 166230.2717381092
$ws.Range("R13").Value2 = 1496072.445642983
$ws.Range("S13").Value2 = 0.1025858189701351
$ws.Range("T13").Value2 = 0.1025858189701351
$ws.Range("G14").Value2 = 62.68962833333333
$ws.Range("H14").Value2 = 188.068885
$ws.Range("I14").Value2 = 0.1415974560293775
$ws.Range("J14").Value2 = 0.1415974560293775
$ws.Range("M14").Value2 = 28.22141
$ws.Range("N14").Value2 = 84.66423
$ws.Range("O14").Value2 = 0.007710741921554872
$ws.Range("P14").Value2 = 0.007710741921554872
$ws.Range("Q14").Value2 = 1769.189703942617
$ws.Range("R14").Value2 = 15922.70733548355
$ws.Range("S14").Value2 = 0.001091821440191243
$ws.Range("T14").Value2 = 0.001091821440191244
$ws.Range("G15").Value2 = 62.68962833333333
$ws.Range("H15").Value2 = 188.068885
$ws.Range("I15").Value2 = 0.1415974560293775
$ws.Range("J15").Value2 = 0.1415974560293775
$ws.Range("O15").Value2 = 0.001803104886918205
$ws.Range("P15").Value2 = 0.001803104886918206
$ws.Range("Q15").Value2 = 413.7130555682033
$ws.Range("R15").Value2 = 3723.41750011383
$ws.Range("S15").Value2 = 0.0002553150649417562
$ws.Range("T15").Value2 = 0.0002553150649417563
$ws.Range("G16").Value2 = 62.68962833333333
$ws.Range("H16").Value2 = 188.068885
$ws.Range("I16").Value2 = 0.1415974560293775
$ws.Range("J16").Value2 = 0.1415974560293775
$ws.Range("M16").Value2 = 1.757996666666666
$ws.Range("N16").Value2 = 5.27399
$ws.Range("O16").Value2 = 0.000480325348578274
$ws.Range("P16").Value2 = 0.0004803253485782741
$ws.Range("Q16").Value2 = 110.2081576445722
$ws.Range("R16").Value2 = 991.8734188011499
$ws.Range("S16").Value2 = 0.00006801284742510756
$ws.Range("T16").Value2 = 0.00006801284742510757
$ws.Range("G17").Value2 = 62.68962833333333
$ws.Range("H17").Value2 = 188.068885
$ws.Range("I17").Value2 = 0.1415974560293775
$ws.Range("J17").Value2 = 0.1415974560293775
$ws.Range("M17").Value2 = 3623.433471666667
$ws.Range("N17").Value2 = 10870.300415
$ws.Range("O17").Value2 = 0.9900058278429487
$ws.Range("P17").Value2 = 0.9900058278429487
$ws.Range("Q17").Value2 = 227151.697629343
$ws.Range("R17").Value2 = 2044365.278664087
$ws.Range("S17").Value2 = 0.1401823066768194
$ws.Range("T17").Value2 = 0.1401823066768194
